# ---------------------------------------------------------------------------
# Edit script: applies the changes described by the commit diff to the
# currently open Word document ($word.ActiveDocument).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for [$findText]"
    }
}

# ---------------------------------------------------------------------------
# 1. Insert a new bold "First Page" paragraph right after the title
#    paragraph, before "Project Name: ...".
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "First Page"
$newPara.Range.Font.Bold = $true

Write-Host "Step 1 done"

# ---------------------------------------------------------------------------
# 2. Date: " [Insert Date]" -> " 2025-02-04"
# ---------------------------------------------------------------------------
Replace-Text "[Insert Date]" "2025-02-04"

Write-Host "Step 2 done"

# ---------------------------------------------------------------------------
# 3. Team member role lines: drop the "proofErr" wrapped " : " typography,
#    add a second role to each member, tweak wording.
# ---------------------------------------------------------------------------
Replace-Text "Full-Stack developer : choose technologies used" `
             "Full-Stack developer-Project Manager: chose technologies used"

Replace-Text "Back-End developer : project definition" `
             "Back-End developer-Database administrator : project definition"

Replace-Text "Front-End developer : Context diagram" `
             "Front-End developer-AI integration specialist : Context diagram"

Replace-Text "Scrum master : UML diagram " `
             "Full-stack developer- Scrum Master : UML diagram "

Write-Host "Step 3 done"
